$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = 44537
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 220
$ws.Range("N9").Value = 3200
$ws.Range("O9").Value = 3500
$ws.Range("P9").Value = 3336
$ws.Range("R9").Value = "Provincia de Linares"
$ws.Range("S9").Value = 1668

$ws.Range("D10").Value = 44516
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 4500
$ws.Range("O10").Value = 5000
$ws.Range("P10").Value = 4750
$ws.Range("R10").Value = "Región de Ñuble"
$ws.Range("S10").Value = 2375

$ws.Range("D11").Value = 44516
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 100
$ws.Range("N11").Value = 4000
$ws.Range("O11").Value = 4000
$ws.Range("P11").Value = 4000
$ws.Range("R11").Value = "Región de Ñuble"
$ws.Range("S11").Value = 2000

$ws.Range("D12").Value = 44187
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 3000
$ws.Range("O12").Value = 3500
$ws.Range("P12").Value = 3250
$ws.Range("R12").Value = "Provincia de Curicó"
$ws.Range("S12").Value = 1625

$ws.Range("D13").Value = 44187
$ws.Range("L13").Value = "Segunda"
$ws.Range("M13").Value = 50
$ws.Range("N13").Value = 2500
$ws.Range("O13").Value = 2500
$ws.Range("P13").Value = 2500
$ws.Range("R13").Value = "Provincia de Curicó"
$ws.Range("S13").Value = 1250

$ws.Range("D14").Value = 44526
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 250
$ws.Range("N14").Value = 3000
$ws.Range("O14").Value = 3500
$ws.Range("P14").Value = 3200
$ws.Range("R14").Value = "Región de Ñuble"
$ws.Range("S14").Value = 1600

$ws.Range("D15").Value = 44523
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 150
$ws.Range("N15").Value = 5000
$ws.Range("O15").Value = 5500
$ws.Range("P15").Value = 5267
$ws.Range("R15").Value = "Provincia de Linares"
$ws.Range("S15").Value = 2634

$ws.Range("D16").Value = 44223
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 3500
$ws.Range("O16").Value = 3500
$ws.Range("P16").Value = 3500
$ws.Range("R16").Value = "Provincia de Curicó"
$ws.Range("S16").Value = 1750

$ws.Range("D17").Value = 44223
$ws.Range("L17").Value = "Segunda"
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 3000
$ws.Range("O17").Value = 3000
$ws.Range("P17").Value = 3000
$ws.Range("R17").Value = "Provincia de Curicó"
$ws.Range("S17").Value = 1500

$ws.Range("D18").Value = 44222
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = 3500
$ws.Range("O18").Value = 4000
$ws.Range("P18").Value = 3750
$ws.Range("R18").Value = "Región de Ñuble"
$ws.Range("S18").Value = 1875

$ws.Range("D19").Value = 44222
$ws.Range("L19").Value = "Segunda"
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = 3000
$ws.Range("O19").Value = 3000
$ws.Range("P19").Value = 3000
$ws.Range("R19").Value = "Región de Ñuble"
$ws.Range("S19").Value = 1500

$ws.Range("D20").Value = 44194
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 400
$ws.Range("N20").Value = 3000
$ws.Range("O20").Value = 3500
$ws.Range("P20").Value = 3250
$ws.Range("R20").Value = "Provincia de Curicó"
$ws.Range("S20").Value = 1625

$ws.Range("D21").Value = 44194
$ws.Range("L21").Value = "Segunda"
$ws.Range("M21").Value = 200
$ws.Range("N21").Value = 2500
$ws.Range("O21").Value = 2500
$ws.Range("P21").Value = 2500
$ws.Range("R21").Value = "Provincia de Curicó"
$ws.Range("S21").Value = 1250

$ws.Range("D22").Value = 44258
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 200
$ws.Range("N22").Value = 3500
$ws.Range("O22").Value = 4000
$ws.Range("P22").Value = 3750
$ws.Range("R22").Value = "Región de Ñuble"
$ws.Range("S22").Value = 1875

$ws.Range("D23").Value = 44258
$ws.Range("L23").Value = "Segunda"
$ws.Range("M23").Value = 100
$ws.Range("N23").Value = 3000
$ws.Range("O23").Value = 3000
$ws.Range("P23").Value = 3000
$ws.Range("R23").Value = "Región de Ñuble"
$ws.Range("S23").Value = 1500

$ws.Range("D24").Value = 44202
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 200
$ws.Range("N24").Value = 2000
$ws.Range("O24").Value = 2500
$ws.Range("P24").Value = 2250
$ws.Range("R24").Value = "Región de O'Higgins"
$ws.Range("S24").Value = 1125

$ws.Range("D25").Value = 44202
$ws.Range("L25").Value = "Segunda"
$ws.Range("M25").Value = 100
$ws.Range("N25").Value = 1500
$ws.Range("O25").Value = 1500
$ws.Range("P25").Value = 1500
$ws.Range("R25").Value = "Región de O'Higgins"
$ws.Range("S25").Value = 750

$ws.Range("D26").Value = 44533
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 200
$ws.Range("N26").Value = 3800
$ws.Range("O26").Value = 4000
$ws.Range("P26").Value = 3900
$ws.Range("R26").Value = "Región de Ñuble"
$ws.Range("S26").Value = 1950

$ws.Range("D27").Value = 44533
$ws.Range("L27").Value = "Segunda"
$ws.Range("M27").Value = 100
$ws.Range("N27").Value = 3500
$ws.Range("O27").Value = 3500
$ws.Range("P27").Value = 3500
$ws.Range("R27").Value = "Región de Ñuble"
$ws.Range("S27").Value = 1750

$ws.Range("D28").Value = 44225
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 100
$ws.Range("N28").Value = 3000
$ws.Range("O28").Value = 3500
$ws.Range("P28").Value = 3250
$ws.Range("R28").Value = "Región de Ñuble"
$ws.Range("S28").Value = 1625

$ws.Range("D29").Value = 44225
$ws.Range("L29").Value = "Segunda"
$ws.Range("M29").Value = 50
$ws.Range("N29").Value = 2500
$ws.Range("O29").Value = 2500
$ws.Range("P29").Value = 2500
$ws.Range("R29").Value = "Región de Ñuble"
$ws.Range("S29").Value = 1250

$ws.Range("D30").Value = 44195
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 100
$ws.Range("N30").Value = 3000
$ws.Range("O30").Value = 3000
$ws.Range("P30").Value = 3000
$ws.Range("R30").Value = "Provincia de Curicó"
$ws.Range("S30").Value = 1500

$ws.Range("D31").Value = 44195
$ws.Range("L31").Value = "Segunda"
$ws.Range("M31").Value = 100
$ws.Range("N31").Value = 2500
$ws.Range("O31").Value = 2500
$ws.Range("P31").Value = 2500
$ws.Range("R31").Value = "Provincia de Curicó"
$ws.Range("S31").Value = 1250

$ws.Range("D32").Value = 44216
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 200
$ws.Range("N32").Value = 3000
$ws.Range("O32").Value = 3500
$ws.Range("P32").Value = 3250
$ws.Range("R32").Value = "Región de Ñuble"
$ws.Range("S32").Value = 1625

$ws.Range("D33").Value = 44216
$ws.Range("L33").Value = "Segunda"
$ws.Range("M33").Value = 100
$ws.Range("N33").Value = 2500
$ws.Range("O33").Value = 2500
$ws.Range("P33").Value = 2500
$ws.Range("R33").Value = "Región de Ñuble"
$ws.Range("S33").Value = 1250

$ws.Range("D34").Value = 44530
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 200
$ws.Range("N34").Value = 5000
$ws.Range("O34").Value = 6000
$ws.Range("P34").Value = 5500
$ws.Range("R34").Value = "Región de Ñuble"
$ws.Range("S34").Value = 2750

$ws.Range("D35").Value = 44530
$ws.Range("L35").Value = "Segunda"
$ws.Range("M35").Value = 100
$ws.Range("N35").Value = 4000
$ws.Range("O35").Value = 4000
$ws.Range("P35").Value = 4000
$ws.Range("R35").Value = "Región de Ñuble"
$ws.Range("S35").Value = 2000

$ws.Range("D36").Value = 44236
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 200
$ws.Range("N36").Value = 3000
$ws.Range("O36").Value = 3500
$ws.Range("P36").Value = 3250
$ws.Range("R36").Value = "Región de Ñuble"
$ws.Range("S36").Value = 1625

$ws.Range("D37").Value = 44236
$ws.Range("L37").Value = "Segunda"
$ws.Range("M37").Value = 100
$ws.Range("N37").Value = 2500
$ws.Range("O37").Value = 2500
$ws.Range("P37").Value = 2500
$ws.Range("R37").Value = "Región de Ñuble"
$ws.Range("S37").Value = 1250

$ws.Range("D38").Value = 44196
$ws.Range("L38").Value = "Primera"
$ws.Range("M38").Value = 200
$ws.Range("N38").Value = 3000
$ws.Range("O38").Value = 3500
$ws.Range("P38").Value = 3250
$ws.Range("R38").Value = "Provincia de Curicó"
$ws.Range("S38").Value = 1625

$ws.Range("D39").Value = 44196
$ws.Range("L39").Value = "Segunda"
$ws.Range("M39").Value = 100
$ws.Range("N39").Value = 2500
$ws.Range("O39").Value = 2500
$ws.Range("P39").Value = 2500
$ws.Range("R39").Value = "Provincia de Curicó"
$ws.Range("S39").Value = 1250

$ws.Range("D40").Value = 44188
$ws.Range("L40").Value = "Primera"
$ws.Range("M40").Value = 200
$ws.Range("N40").Value = 3000
$ws.Range("O40").Value = 3500
$ws.Range("P40").Value = 3250
$ws.Range("R40").Value = "Provincia de Curicó"
$ws.Range("S40").Value = 1625

$ws.Range("D41").Value = 44159
$ws.Range("L41").Value = "Primera"
$ws.Range("M41").Value = 100
$ws.Range("N41").Value = 6000
$ws.Range("O41").Value = 6000
$ws.Range("P41").Value = 6000
$ws.Range("R41").Value = "Provincia de Curicó"
$ws.Range("S41").Value = 3000

$ws.Range("D42").Value = 44159
$ws.Range("L42").Value = "Segunda"
$ws.Range("M42").Value = 100
$ws.Range("N42").Value = 5000
$ws.Range("O42").Value = 5000
$ws.Range("P42").Value = 5000
$ws.Range("R42").Value = "Provincia de Curicó"
$ws.Range("S42").Value = 2500

$ws.Range("D43").Value = 44505
$ws.Range("L43").Value = "Segunda"
$ws.Range("M43").Value = 100
$ws.Range("N43").Value = 7000
$ws.Range("O43").Value = 7000
$ws.Range("P43").Value = 7000
$ws.Range("R43").Value = "Provincia de Curicó"
$ws.Range("S43").Value = 3500

$ws.Range("D44").Value = 44208
$ws.Range("L44").Value = "Primera"
$ws.Range("M44").Value = 100
$ws.Range("N44").Value = 3000
$ws.Range("O44").Value = 3500
$ws.Range("P44").Value = 3250
$ws.Range("R44").Value = "Región de Ñuble"
$ws.Range("S44").Value = 1625

$ws.Range("D45").Value = 44208
$ws.Range("L45").Value = "Segunda"
$ws.Range("M45").Value = 50
$ws.Range("N45").Value = 2500
$ws.Range("O45").Value = 2500
$ws.Range("P45").Value = 2500
$ws.Range("R45").Value = "Región de Ñuble"
$ws.Range("S45").Value = 1250
